$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Reword the results paragraph describing the climbing assay.
#    Locate the paragraph by its stable leading text, then replace the
#    trailing portion (everything after "...flies that eclose") with the
#    new wording. The italic "Raptor" run that used to sit in the middle of
#    this passage is removed entirely, so the whole replaced span becomes
#    plain (non-italic) text, matching the body font already used there.
# ---------------------------------------------------------------------------

$anchorText = "suppression leads to reduced muscle function in the flies that eclose"

$find = $d.Content
$found = $find.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor text for the climbing-assay paragraph."
}

$anchorEnd = $find.End

# Figure out where this run-on sentence actually ends (the paragraph mark),
# by walking forward to the paragraph containing the anchor.
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $find.Start -and $find.End -le $p.Range.End) {
        $para = $p
        break
    }
}
if ($para -eq $null) {
    throw "Could not locate paragraph for the climbing-assay text."
}

# End of the paragraph's visible text, i.e. right before the paragraph mark.
$paraTextEnd = $para.Range.End - 1

$tail = $d.Range($anchorEnd, $paraTextEnd)

$newTail = " even very early, consistent a developmental problem in myogenesis" + `
    ".  Interestingly, " + `
    "these problems persist throughout the lifespan of the fly, even in those that reach adulthood.  Also interesting, is that there was a correspondence between the efficienty of the shRNA strain to cause lethality and its effects on muscle function, indicating a potential gene-dosage effect on both of these phenotypes."

$tail.Text = $newTail

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the Figure 5 caption (after "(in
#    days)") to the end of the paragraph we just edited. Word keeps only a
#    single "_GoBack" bookmark, so adding it at the new spot implicitly
#    removes it from the old one.
# ---------------------------------------------------------------------------

# NOTE: adding a bookmark whose zero-length range sits exactly at a
# paragraph's end position (i.e. immediately before the paragraph mark)
# confuses the engine and the bookmark ends up at the very start of the
# document instead. Work around it by temporarily inserting a sentinel
# character after the target spot, anchoring the bookmark next to that
# (now interior) position, and then removing the sentinel again.
$newBookmarkPos = $anchorEnd + $newTail.Length
$sentinelPt = $d.Range($newBookmarkPos, $newBookmarkPos)
$sentinelPt.InsertAfter("Z")

$bookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$sentinelRange = $d.Range($newBookmarkPos, $newBookmarkPos + 1)
$sentinelRange.Delete()

Write-Host "Done."
